# "Updated comparison 9 channels by new results of VSC"
#
# The SN9_RunTime sheet gets a new row inserted between the existing
# "Dekempeneer & Derkinderen (20 threads)" row (row 2) and the
# "Codish et al. (288 threads)" row (old row 3, becomes row 4):
# a new "Algoritme Codish et al. door / Dekempeneer & Derkinderen /
# (20 threads)" row with its own runtime figures (in minutes/seconds,
# derived from a raw nanosecond measurement).
#
# Inserting the row natively keeps styles + shifts/repairs every
# dependent formula (e.g. F2's D2/D3 -> D2/D4) automatically, exactly
# like using Excel's own "Insert Sheet Rows" on row 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SN9_RunTime")

# Remember whichever sheet is active right now so we don't change the
# workbook's active tab just to set SN9_RunTime's own cell selection.
$previouslyActive = $wb.ActiveSheet

# Push the old row 3 ("Codish et al.") down to row 4 and open up a
# fresh (style-inheriting) row 3.
$ws.Rows.Item(3).Insert()

# New label for A3.
$ws.Range("A3").Value = "Algoritme Codish et al. door`nDekempeneer & Derkinderen`n(20 threads)"

# Raw measurement (ns) -> minutes, mirroring the existing B/C/D column
# pattern (D = raw converted to seconds-ish base unit, C = D/60, B = C/60).
$ws.Range("D3").Formula = "=22308959989475/1000000000"
$ws.Range("C3").Formula = "=D3/60"
$ws.Range("B3").Formula = "=C3/60"

# F3/G3 stay present (same style as F2/G2) but empty, matching the old
# row 3's F/G cells before the insert.
$ws.Range("F3").Value = ""
$ws.Range("G3").Value = ""

# The new row is taller (wraps 3 lines of text) than the default/auto
# height Excel would otherwise give it.
$ws.Rows.Item(3).RowHeight = 75

# Put the selection on B3 (as in the saved file) without leaving
# SN9_RunTime as the active/selected sheet tab.
$ws.Activate()
$ws.Range("B3").Select()
$previouslyActive.Activate()
